# "backend works on azure now"
#
# The canonical-XML diff for this commit touches only two things:
#
#   1. ppt/presentation.xml (+ the slide's <we:webextensionref r:id=...> and
#      the fallback <a:blip r:embed=...>): every r:id="R<32 hex chars>"
#      token is replaced by a *different* random-looking r:id. These are
#      OPC relationship-id strings that PowerPoint mints fresh, arbitrarily,
#      every time it rewrites a package; the targets they point at
#      (slideMaster -> slideMaster1.xml, slide -> slide1.xml, the image
#      relationship, the webextension relationship, ...) are completely
#      unchanged before/after. There is no content change to replay here --
#      it is just resave churn, and it is not something any host-level
#      automation call (PowerPoint COM/VBA included) can target, since
#      relationship ids are never surfaced as settable values on the object
#      model; the host assigns them internally.
#
#   2. ppt/slides/udata/data.xml: the <we:webextension id="{...guid...}">
#      attribute changes to a new GUID (the add-in instance id PowerPoint
#      recorded when the task-pane add-in -- the "OfficeApp" content -- was
#      last (re)inserted/reconnected, consistent with the add-in now
#      talking to its Azure-hosted backend). This value lives entirely
#      inside the webextension part that backs the slide's "OfficeApp 0"
#      AlternateContent shape; the PowerPoint object model does not expose
#      webextension parts (task-pane add-in bindings/snapshots) for editing
#      at all -- there is no Shape/Presentation property or method for it,
#      by design (add-ins are inserted/managed through the Office Store /
#      "Insert Add-in" UI, never through Shapes/TextFrame-style automation).
#      Confirmed empirically against this host too: Shapes.Item/Range can't
#      even address that graphicFrame distinctly from "Title 1" (both
#      happen to share cNvPr id="2" in this deck), so any attempted write
#      through the Shapes collection lands on the title placeholder instead
#      -- i.e. touching it here would actively corrupt shape 1 rather than
#      reach the intended part.
#
# Net effect: nothing in this commit is reachable/safe to perform through
# PowerPoint COM automation (no slide text, shape geometry, or visible
# content changed). So this script intentionally makes no Shapes/TextFrame
# edits -- doing so would only risk corrupting the Title/Subtitle
# placeholders for no real gain -- and simply round-trips the deck, which
# is the faithful behavior for an edit whose only substance is an internal
# add-in identifier minted by PowerPoint itself.

$p = $ppt.ActivePresentation

# Touch nothing on the slide (Title/Subtitle text and the OfficeApp
# AlternateContent block are all unchanged in the diff); just confirm the
# deck is in the expected shape and let the host persist it as-is.
$s = $p.Slides.Item(1)
Write-Host ("Slides: {0}, Shapes on slide 1: {1}" -f $p.Slides.Count, $s.Shapes.Count)

$p.Save()
